# Automatic map update (2025-07-29 09:07:54)
# Inserts a new incident row at row 75, shifting the previously existing
# rows 75-79 down to rows 76-80.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above current row 75, shifting rows 75-79 (and
# everything below them) down by one row.
$ws.Rows.Item(75).Insert()

# Columns A, B, D and E hold values that look numeric/date-like
# ("6522", "7/28/2025", "7", "808571978") but must stay plain text, just
# like every other row in this sheet. Temporarily force a text number
# format so Excel doesn't auto-convert them, then clear the (no longer
# needed) formatting once the literal text value has been stored.
$textCols = $ws.Range("A75:E75")
$textCols.NumberFormat = "@"

# Populate the newly inserted row 75 with the new record's data.
$ws.Cells.Item(75, 1).Value = "6522"
$ws.Cells.Item(75, 2).Value = "7/28/2025"
$ws.Cells.Item(75, 3).Value = "ESTADO PLURINACIONAL DE BOLIVIA 384"
$ws.Cells.Item(75, 4).Value = "7"
$ws.Cells.Item(75, 5).Value = "808571978"
$ws.Cells.Item(75, 6).Value = "PEBCOM"
$ws.Cells.Item(75, 7).Value = "Pendiente"
$ws.Cells.Item(75, 8).Value = "Traspasar red y desmontar"
$ws.Cells.Item(75, 9).Value = 1
$ws.Cells.Item(75, 10).Value = "Cambio"
$ws.Cells.Item(75, 11).Value = "Sin equipos"
$ws.Cells.Item(75, 12).Value = "Pasante"
$ws.Cells.Item(75, 13).Value = -58.466995
$ws.Cells.Item(75, 14).Value = -34.626426
$ws.Cells.Item(75, 15).Value = "Boedo"
$ws.Cells.Item(75, 16).Value = "Capital Sur"

$textCols.ClearFormats()
